$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: JACQUELINE MARTINEZ ARELLANO, doc 45548983, periodo 1905, valor mora 33125, salario 1300000
$ws.Range("C16").Value = "45548983"
$ws.Range("D16").Value = "JACQUELINE MARTINEZ ARELLANO"
$ws.Range("E16").Value = "1905"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 1300000

# Row 17: JACQUELINE MARTINEZ ARELLANO, doc 45548983, periodo 1906, valor mora 33125, salario 1300000
$ws.Range("C17").Value = "45548983"
$ws.Range("D17").Value = "JACQUELINE MARTINEZ ARELLANO"
$ws.Range("E17").Value = "1906"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 1300000

# Row 18: YULIS ROXANA MARTINEZ GONZALEZ, doc 1002241233, periodo 2001, valor mora 39227, salario 1000000
$ws.Range("C18").Value = "1002241233"
$ws.Range("D18").Value = "YULIS ROXANA MARTINEZ GONZALEZ"
$ws.Range("E18").Value = "2001"
$ws.Range("F18").Value = 39227
$ws.Range("G18").Value = 1000000
